# Auto update Excel log
# Appends new sensor/log rows to three sheets: ALERTS, Proximity, Camera.
$wb = $excel.ActiveWorkbook

# --- ALERTS sheet: two new FALL_DETECTED rows (15 and 16) ---
$ws = $wb.Worksheets.Item("ALERTS")

$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "2026-02-01"
$ws.Range("B15").Value = "14:43:25"
$ws.Range("C15").Value = "14:00"
$ws.Range("D15").Value = "Living Room"
$ws.Range("E15").Value = "CRITICAL"
$ws.Range("F15").Value = "FALL_DETECTED"

$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "2026-02-01"
$ws.Range("B16").Value = "14:43:26"
$ws.Range("C16").Value = "14:00"
$ws.Range("D16").Value = "Living Room"
$ws.Range("E16").Value = "CRITICAL"
$ws.Range("F16").Value = "FALL_DETECTED"

# --- Proximity sheet: one new ENTER row (39) ---
$ws = $wb.Worksheets.Item("Proximity")

$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "2026-02-01"
$ws.Range("B39").Value = "14:43:34"
$ws.Range("C39").Value = "14:00"
$ws.Range("D39").Value = "Living Room Main Door"
$ws.Range("E39").Value = "ENTER"
$ws.Range("F39").Value = "User ENTERED Living Room Main Door"

# --- Camera sheet: one new Image Captured row (25) ---
$ws = $wb.Worksheets.Item("Camera")

$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "2026-02-01"
$ws.Range("B25").Value = "14:43:35"
$ws.Range("C25").Value = "14:00"
$ws.Range("D25").Value = "Living Room Main Door"
$ws.Range("E25").Value = "Image Captured"
$ws.Range("F25").Value = "Active"

Write-Output "Added rows: ALERTS 15-16, Proximity 39, Camera 25"
